# Applies cached-value corrections to the Leve profit/price columns (H:N)
# across all 8 job sheets, per the upstream data refresh.
# Values are plain numbers (no formulas are used in this workbook).
$wb = $excel.ActiveWorkbook


# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
# row 2
$ws.Cells.Item(2,8).Value = 111.57143  # H2
$ws.Cells.Item(2,9).Value = 111.57143  # I2
$ws.Cells.Item(2,11).Value = 111.57143  # K2
$ws.Cells.Item(2,13).Value = 1.428569999999993  # M2
# row 11
$ws.Cells.Item(11,8).Value = 16.75  # H11
$ws.Cells.Item(11,9).Value = 16.75  # I11
$ws.Cells.Item(11,11).Value = 16.75  # K11
$ws.Cells.Item(11,13).Value = 123.25  # M11
# row 19
$ws.Cells.Item(19,8).Value = 799.2381  # H19
$ws.Cells.Item(19,9).Value = 1202  # I19
$ws.Cells.Item(19,10).Value = 597.8570999999999  # J19
$ws.Cells.Item(19,11).Value = 1202  # K19
$ws.Cells.Item(19,12).Value = 597.8570999999999  # L19
$ws.Cells.Item(19,13).Value = -1027  # M19
$ws.Cells.Item(19,14).Value = -947.8570999999999  # N19
# row 28
$ws.Cells.Item(28,8).Value = 3138.1428  # H28
$ws.Cells.Item(28,9).Value = 1567.7273  # I28
$ws.Cells.Item(28,10).Value = 8896.333000000001  # J28
$ws.Cells.Item(28,11).Value = 1567.7273  # K28
$ws.Cells.Item(28,12).Value = 8896.333000000001  # L28
$ws.Cells.Item(28,13).Value = -1082.7273  # M28
$ws.Cells.Item(28,14).Value = -9866.333000000001  # N28
# row 32
$ws.Cells.Item(32,8).Value = 749.5  # H32
$ws.Cells.Item(32,9).Value = 500  # I32
$ws.Cells.Item(32,11).Value = 500  # K32
$ws.Cells.Item(32,13).Value = -174  # M32
# row 38
$ws.Cells.Item(38,8).Value = 2731.077  # H38
$ws.Cells.Item(38,10).Value = 5313.6  # J38
$ws.Cells.Item(38,12).Value = 15940.8  # L38
$ws.Cells.Item(38,14).Value = -16684.8  # N38
# row 39
$ws.Cells.Item(39,8).Value = 260  # H39
$ws.Cells.Item(39,9).Value = 260  # I39
$ws.Cells.Item(39,10).Value = 0  # J39
$ws.Cells.Item(39,11).Value = 780  # K39
$ws.Cells.Item(39,12).Value = 0  # L39
$ws.Cells.Item(39,13).Value = -484  # M39
$ws.Cells.Item(39,14).ClearContents()  # N39 removed
# row 41
$ws.Cells.Item(41,8).Value = 1484.6923  # H41
$ws.Cells.Item(41,9).Value = 412.25  # I41
$ws.Cells.Item(41,10).Value = 1961.3334  # J41
$ws.Cells.Item(41,11).Value = 412.25  # K41
$ws.Cells.Item(41,12).Value = 1961.3334  # L41
$ws.Cells.Item(41,13).Value = 27.75  # M41
$ws.Cells.Item(41,14).Value = -2841.3334  # N41
# row 43
$ws.Cells.Item(43,8).Value = 10012049  # H43
$ws.Cells.Item(43,9).Value = 14287928  # I43
$ws.Cells.Item(43,10).Value = 35000  # J43
$ws.Cells.Item(43,11).Value = 14287928  # K43
$ws.Cells.Item(43,12).Value = 35000  # L43
$ws.Cells.Item(43,13).Value = -14287859  # M43
$ws.Cells.Item(43,14).Value = -35138  # N43
# row 45
$ws.Cells.Item(45,8).Value = 1436  # H45
$ws.Cells.Item(45,9).Value = 1100  # I45
$ws.Cells.Item(45,10).Value = 1484  # J45
$ws.Cells.Item(45,11).Value = 3300  # K45
$ws.Cells.Item(45,12).Value = 4452  # L45
$ws.Cells.Item(45,13).Value = -3108  # M45
$ws.Cells.Item(45,14).Value = -4836  # N45
# row 53
$ws.Cells.Item(53,8).Value = 322.33334  # H53
$ws.Cells.Item(53,9).Value = 257.2857  # I53
$ws.Cells.Item(53,11).Value = 257.2857  # K53
$ws.Cells.Item(53,13).Value = 379.7143  # M53
# row 61
$ws.Cells.Item(61,8).Value = 265  # H61
$ws.Cells.Item(61,9).Value = 265  # I61
$ws.Cells.Item(61,11).Value = 795  # K61
$ws.Cells.Item(61,13).Value = -623  # M61
# row 62
$ws.Cells.Item(62,8).Value = 9628.091  # H62
$ws.Cells.Item(62,9).Value = 7651.5  # I62
$ws.Cells.Item(62,11).Value = 7651.5  # K62
$ws.Cells.Item(62,13).Value = -7027.5  # M62
# row 65
$ws.Cells.Item(65,8).Value = 9628.091  # H65
$ws.Cells.Item(65,9).Value = 7651.5  # I65
$ws.Cells.Item(65,11).Value = 38257.5  # K65
$ws.Cells.Item(65,13).Value = -35137.5  # M65
# row 69
$ws.Cells.Item(69,8).Value = 7069.125  # H69
$ws.Cells.Item(69,10).Value = 7069.125  # J69
$ws.Cells.Item(69,12).Value = 21207.375  # L69
$ws.Cells.Item(69,14).Value = -22955.375  # N69
# row 72
$ws.Cells.Item(72,8).Value = 7069.125  # H72
$ws.Cells.Item(72,10).Value = 7069.125  # J72
$ws.Cells.Item(72,12).Value = 63622.125  # L72
$ws.Cells.Item(72,14).Value = -72358.125  # N72
# row 76
$ws.Cells.Item(76,8).Value = 0  # H76
$ws.Cells.Item(76,9).Value = 0  # I76
$ws.Cells.Item(76,11).Value = 0  # K76
$ws.Cells.Item(76,13).ClearContents()  # M76 removed
# row 79
$ws.Cells.Item(79,8).Value = 0  # H79
$ws.Cells.Item(79,9).Value = 0  # I79
$ws.Cells.Item(79,11).Value = 0  # K79
$ws.Cells.Item(79,13).ClearContents()  # M79 removed
# row 87
$ws.Cells.Item(87,8).Value = 149420  # H87
$ws.Cells.Item(87,10).Value = 149420  # J87
$ws.Cells.Item(87,12).Value = 149420  # L87
$ws.Cells.Item(87,14).Value = -151916  # N87
# row 88
$ws.Cells.Item(88,8).Value = 833  # H88
$ws.Cells.Item(88,10).Value = 899.5  # J88
$ws.Cells.Item(88,12).Value = 899.5  # L88
$ws.Cells.Item(88,14).Value = -1711.5  # N88
# row 90
$ws.Cells.Item(90,8).Value = 149420  # H90
$ws.Cells.Item(90,10).Value = 149420  # J90
$ws.Cells.Item(90,12).Value = 448260  # L90
$ws.Cells.Item(90,14).Value = -460740  # N90
# row 91
$ws.Cells.Item(91,8).Value = 833  # H91
$ws.Cells.Item(91,10).Value = 899.5  # J91
$ws.Cells.Item(91,12).Value = 899.5  # L91
$ws.Cells.Item(91,14).Value = -3707.5  # N91
# row 92
$ws.Cells.Item(92,8).Value = 1703.4736  # H92
$ws.Cells.Item(92,9).Value = 1850.1333  # I92
$ws.Cells.Item(92,10).Value = 1153.5  # J92
$ws.Cells.Item(92,11).Value = 1850.1333  # K92
$ws.Cells.Item(92,12).Value = 1153.5  # L92
$ws.Cells.Item(92,13).Value = -602.1333  # M92
$ws.Cells.Item(92,14).Value = -3649.5  # N92
# row 98
$ws.Cells.Item(98,8).Value = 1167.3334  # H98
$ws.Cells.Item(98,9).Value = 1167.3334  # I98
$ws.Cells.Item(98,11).Value = 1167.3334  # K98
$ws.Cells.Item(98,13).Value = 330.6666  # M98
# row 99
$ws.Cells.Item(99,8).Value = 2908.7778  # H99
$ws.Cells.Item(99,9).Value = 250  # I99
$ws.Cells.Item(99,10).Value = 4238.1665  # J99
$ws.Cells.Item(99,11).Value = 750  # K99
$ws.Cells.Item(99,12).Value = 12714.4995  # L99
$ws.Cells.Item(99,13).Value = 748  # M99
$ws.Cells.Item(99,14).Value = -15710.4995  # N99
# row 101
$ws.Cells.Item(101,8).Value = 787  # H101
$ws.Cells.Item(101,9).Value = 724.5  # I101
$ws.Cells.Item(101,10).Value = 849.5  # J101
$ws.Cells.Item(101,11).Value = 2173.5  # K101
$ws.Cells.Item(101,12).Value = 2548.5  # L101
$ws.Cells.Item(101,13).Value = -551.5  # M101
$ws.Cells.Item(101,14).Value = -5792.5  # N101
# row 103
$ws.Cells.Item(103,8).Value = 1850.3684  # H103
$ws.Cells.Item(103,9).Value = 2828.5  # I103
$ws.Cells.Item(103,10).Value = 1398.9231  # J103
$ws.Cells.Item(103,11).Value = 8485.5  # K103
$ws.Cells.Item(103,12).Value = 4196.7693  # L103
$ws.Cells.Item(103,13).Value = -7899.5  # M103
$ws.Cells.Item(103,14).Value = -5368.7693  # N103
# row 122
$ws.Cells.Item(122,8).Value = 1167.3334  # H122
$ws.Cells.Item(122,9).Value = 1167.3334  # I122
$ws.Cells.Item(122,11).Value = 3502.0002  # K122
$ws.Cells.Item(122,13).Value = -1052.0002  # M122
# row 129
$ws.Cells.Item(129,8).Value = 1404.6666  # H129
$ws.Cells.Item(129,9).Value = 685.6  # I129
$ws.Cells.Item(129,10).Value = 5000  # J129
$ws.Cells.Item(129,11).Value = 2056.8  # K129
$ws.Cells.Item(129,12).Value = 15000  # L129
$ws.Cells.Item(129,13).Value = 2943.2  # M129
$ws.Cells.Item(129,14).Value = -25000  # N129
# row 132
$ws.Cells.Item(132,8).Value = 27894.125  # H132
$ws.Cells.Item(132,9).Value = 43077.6  # I132
$ws.Cells.Item(132,11).Value = 129232.8  # K132
$ws.Cells.Item(132,13).Value = -126702.8  # M132
# row 137
$ws.Cells.Item(137,8).Value = 2902.0952  # H137
$ws.Cells.Item(137,9).Value = 1291.3334  # I137
$ws.Cells.Item(137,10).Value = 3170.5557  # J137
$ws.Cells.Item(137,11).Value = 3874.0002  # K137
$ws.Cells.Item(137,12).Value = 9511.667099999999  # L137
$ws.Cells.Item(137,13).Value = -1324.0002  # M137
$ws.Cells.Item(137,14).Value = -14611.6671  # N137
# row 138
$ws.Cells.Item(138,8).Value = 4689.3125  # H138
$ws.Cells.Item(138,10).Value = 4719.5835  # J138
$ws.Cells.Item(138,12).Value = 14158.7505  # L138
$ws.Cells.Item(138,14).Value = -24438.7505  # N138
# row 141
$ws.Cells.Item(141,8).Value = 4246.25  # H141
$ws.Cells.Item(141,9).Value = 5333.3335  # I141
$ws.Cells.Item(141,10).Value = 985  # J141
$ws.Cells.Item(141,11).Value = 16000.0005  # K141
$ws.Cells.Item(141,12).Value = 2955  # L141
$ws.Cells.Item(141,13).Value = -10820.0005  # M141
$ws.Cells.Item(141,14).Value = -13315  # N141

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
# row 2
$ws.Cells.Item(2,8).Value = 681.94446  # H2
$ws.Cells.Item(2,9).Value = 569.9286  # I2
$ws.Cells.Item(2,10).Value = 1074  # J2
$ws.Cells.Item(2,11).Value = 569.9286  # K2
$ws.Cells.Item(2,12).Value = 1074  # L2
$ws.Cells.Item(2,13).Value = -456.9286  # M2
$ws.Cells.Item(2,14).Value = -1300  # N2
# row 5
$ws.Cells.Item(5,8).Value = 857.1429000000001  # H5
$ws.Cells.Item(5,9).Value = 857.1429000000001  # I5
$ws.Cells.Item(5,11).Value = 857.1429000000001  # K5
$ws.Cells.Item(5,13).Value = -745.1429000000001  # M5
# row 32
$ws.Cells.Item(32,8).Value = 2176412.2  # H32
$ws.Cells.Item(32,9).Value = 2307.5527  # I32
$ws.Cells.Item(32,10).Value = 12503410  # J32
$ws.Cells.Item(32,11).Value = 2307.5527  # K32
$ws.Cells.Item(32,12).Value = 12503410  # L32
$ws.Cells.Item(32,13).Value = -2020.5527  # M32
$ws.Cells.Item(32,14).Value = -12503984  # N32
# row 45
$ws.Cells.Item(45,8).Value = 2938  # H45
$ws.Cells.Item(45,9).Value = 2394.5293  # I45
$ws.Cells.Item(45,11).Value = 2394.5293  # K45
$ws.Cells.Item(45,13).Value = -2017.5293  # M45
# row 97
$ws.Cells.Item(97,8).Value = 3000  # H97
$ws.Cells.Item(97,9).Value = 3000  # I97
$ws.Cells.Item(97,10).Value = 3000  # J97
$ws.Cells.Item(97,11).Value = 3000  # K97
$ws.Cells.Item(97,12).Value = 3000  # L97
$ws.Cells.Item(97,13).Value = -2504  # M97
$ws.Cells.Item(97,14).Value = -3992  # N97
# row 116
$ws.Cells.Item(116,8).Value = 681.94446  # H116
$ws.Cells.Item(116,9).Value = 569.9286  # I116
$ws.Cells.Item(116,10).Value = 1074  # J116
$ws.Cells.Item(116,11).Value = 569.9286  # K116
$ws.Cells.Item(116,12).Value = 1074  # L116
$ws.Cells.Item(116,13).Value = 1724.0714  # M116
$ws.Cells.Item(116,14).Value = -5662  # N116
# row 122
$ws.Cells.Item(122,8).Value = 4413.25  # H122
$ws.Cells.Item(122,9).Value = 4413.25  # I122
$ws.Cells.Item(122,10).Value = 0  # J122
$ws.Cells.Item(122,11).Value = 13239.75  # K122
$ws.Cells.Item(122,12).Value = 0  # L122
$ws.Cells.Item(122,13).Value = -10789.75  # M122
$ws.Cells.Item(122,14).ClearContents()  # N122 removed
# row 131
$ws.Cells.Item(131,8).Value = 70000  # H131
$ws.Cells.Item(131,10).Value = 70000  # J131
$ws.Cells.Item(131,12).Value = 70000  # L131
$ws.Cells.Item(131,14).Value = -80080  # N131

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
# row 3
$ws.Cells.Item(3,8).Value = 681.94446  # H3
$ws.Cells.Item(3,9).Value = 569.9286  # I3
$ws.Cells.Item(3,10).Value = 1074  # J3
$ws.Cells.Item(3,11).Value = 569.9286  # K3
$ws.Cells.Item(3,12).Value = 1074  # L3
$ws.Cells.Item(3,13).Value = -455.9286  # M3
$ws.Cells.Item(3,14).Value = -1302  # N3
# row 4
$ws.Cells.Item(4,8).Value = 857.1429000000001  # H4
$ws.Cells.Item(4,9).Value = 857.1429000000001  # I4
$ws.Cells.Item(4,11).Value = 857.1429000000001  # K4
$ws.Cells.Item(4,13).Value = -742.1429000000001  # M4
# row 20
$ws.Cells.Item(20,8).Value = 2351.75  # H20
$ws.Cells.Item(20,9).Value = 1899  # I20
$ws.Cells.Item(20,11).Value = 1899  # K20
$ws.Cells.Item(20,13).Value = -1652  # M20
# row 22
$ws.Cells.Item(22,8).Value = 1559.5385  # H22
$ws.Cells.Item(22,9).Value = 1479.4546  # I22
$ws.Cells.Item(22,11).Value = 1479.4546  # K22
$ws.Cells.Item(22,13).Value = -1306.4546  # M22
# row 80
$ws.Cells.Item(80,8).Value = 580.6429000000001  # H80
$ws.Cells.Item(80,9).Value = 367.25  # I80
$ws.Cells.Item(80,11).Value = 367.25  # K80
$ws.Cells.Item(80,13).Value = 630.75  # M80
# row 83
$ws.Cells.Item(83,8).Value = 580.6429000000001  # H83
$ws.Cells.Item(83,9).Value = 367.25  # I83
$ws.Cells.Item(83,11).Value = 1836.25  # K83
$ws.Cells.Item(83,13).Value = 3155.75  # M83
# row 94
$ws.Cells.Item(94,8).Value = 1649.1428  # H94
$ws.Cells.Item(94,10).Value = 0  # J94
$ws.Cells.Item(94,12).Value = 0  # L94
$ws.Cells.Item(94,14).ClearContents()  # N94 removed
# row 102
$ws.Cells.Item(102,8).Value = 15556  # H102
$ws.Cells.Item(102,9).Value = 15556  # I102
$ws.Cells.Item(102,11).Value = 15556  # K102
$ws.Cells.Item(102,13).Value = -12311  # M102
# row 105
$ws.Cells.Item(105,8).Value = 1950  # H105
$ws.Cells.Item(105,9).Value = 1950  # I105
$ws.Cells.Item(105,11).Value = 1950  # K105
$ws.Cells.Item(105,13).Value = -203  # M105
# row 107
$ws.Cells.Item(107,8).Value = 5722  # H107
$ws.Cells.Item(107,9).Value = 3583  # I107
$ws.Cells.Item(107,11).Value = 3583  # K107
$ws.Cells.Item(107,13).Value = -1663  # M107
# row 134
$ws.Cells.Item(134,8).Value = 3955.8333  # H134
$ws.Cells.Item(134,9).Value = 4082  # I134
$ws.Cells.Item(134,10).Value = 3325  # J134
$ws.Cells.Item(134,11).Value = 12246  # K134
$ws.Cells.Item(134,12).Value = 9975  # L134
$ws.Cells.Item(134,13).Value = -9711  # M134
$ws.Cells.Item(134,14).Value = -15045  # N134

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
# row 6
$ws.Cells.Item(6,8).Value = 2034  # H6
$ws.Cells.Item(6,9).Value = 2034  # I6
$ws.Cells.Item(6,11).Value = 2034  # K6
$ws.Cells.Item(6,13).Value = -1921  # M6
# row 10
$ws.Cells.Item(10,8).Value = 2375  # H10
$ws.Cells.Item(10,9).Value = 2375  # I10
$ws.Cells.Item(10,11).Value = 2375  # K10
$ws.Cells.Item(10,13).Value = -2236  # M10
# row 31
$ws.Cells.Item(31,8).Value = 4292.75  # H31
$ws.Cells.Item(31,9).Value = 2427.9285  # I31
$ws.Cells.Item(31,10).Value = 4860.304  # J31
$ws.Cells.Item(31,11).Value = 2427.9285  # K31
$ws.Cells.Item(31,12).Value = 4860.304  # L31
$ws.Cells.Item(31,13).Value = -2132.9285  # M31
$ws.Cells.Item(31,14).Value = -5450.304  # N31
# row 33
$ws.Cells.Item(33,8).Value = 15260.125  # H33
$ws.Cells.Item(33,9).Value = 15260.125  # I33
$ws.Cells.Item(33,11).Value = 15260.125  # K33
$ws.Cells.Item(33,13).Value = -14881.125  # M33
# row 34
$ws.Cells.Item(34,8).Value = 4292.75  # H34
$ws.Cells.Item(34,9).Value = 2427.9285  # I34
$ws.Cells.Item(34,10).Value = 4860.304  # J34
$ws.Cells.Item(34,11).Value = 2427.9285  # K34
$ws.Cells.Item(34,12).Value = 4860.304  # L34
$ws.Cells.Item(34,13).Value = -2225.9285  # M34
$ws.Cells.Item(34,14).Value = -5264.304  # N34
# row 59
$ws.Cells.Item(59,8).Value = 141187200  # H59
$ws.Cells.Item(59,10).Value = 141187200  # J59
$ws.Cells.Item(59,12).Value = 141187200  # L59
$ws.Cells.Item(59,14).Value = -141189490  # N59
# row 94
$ws.Cells.Item(94,8).Value = 3369  # H94
$ws.Cells.Item(94,9).Value = 1556.1111  # I94
$ws.Cells.Item(94,10).Value = 6632.2  # J94
$ws.Cells.Item(94,11).Value = 1556.1111  # K94
$ws.Cells.Item(94,12).Value = 6632.2  # L94
$ws.Cells.Item(94,13).Value = -1105.1111  # M94
$ws.Cells.Item(94,14).Value = -7534.2  # N94
# row 105
$ws.Cells.Item(105,8).Value = 1336.4762  # H105
$ws.Cells.Item(105,9).Value = 798.38464  # I105
$ws.Cells.Item(105,10).Value = 2210.875  # J105
$ws.Cells.Item(105,11).Value = 798.38464  # K105
$ws.Cells.Item(105,12).Value = 2210.875  # L105
$ws.Cells.Item(105,13).Value = 948.61536  # M105
$ws.Cells.Item(105,14).Value = -5704.875  # N105
# row 130
$ws.Cells.Item(130,8).Value = 26249.75  # H130
$ws.Cells.Item(130,10).Value = 26249.75  # J130
$ws.Cells.Item(130,12).Value = 26249.75  # L130
$ws.Cells.Item(130,14).Value = -36289.75  # N130
# row 132
$ws.Cells.Item(132,8).Value = 1724.3636  # H132
$ws.Cells.Item(132,9).Value = 1142  # I132
$ws.Cells.Item(132,11).Value = 3426  # K132
$ws.Cells.Item(132,13).Value = -896  # M132
# row 134
$ws.Cells.Item(134,8).Value = 1082  # H134
$ws.Cells.Item(134,9).Value = 1085.7273  # I134
$ws.Cells.Item(134,10).Value = 1000  # J134
$ws.Cells.Item(134,11).Value = 3257.1819  # K134
$ws.Cells.Item(134,12).Value = 3000  # L134
$ws.Cells.Item(134,13).Value = -722.1819  # M134
$ws.Cells.Item(134,14).Value = -8070  # N134

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
# row 7
$ws.Cells.Item(7,8).Value = 201.66667  # H7
$ws.Cells.Item(7,10).Value = 420.3  # J7
$ws.Cells.Item(7,12).Value = 1260.9  # L7
$ws.Cells.Item(7,14).Value = -1484.9  # N7
# row 16
$ws.Cells.Item(16,8).Value = 440  # H16
$ws.Cells.Item(16,10).Value = 440  # J16
$ws.Cells.Item(16,12).Value = 1320  # L16
$ws.Cells.Item(16,14).Value = -1666  # N16
# row 23
$ws.Cells.Item(23,8).Value = 104.42857  # H23
$ws.Cells.Item(23,10).Value = 106.2  # J23
$ws.Cells.Item(23,12).Value = 318.6  # L23
$ws.Cells.Item(23,14).Value = -788.6  # N23
# row 46
$ws.Cells.Item(46,8).Value = 1999  # H46
$ws.Cells.Item(46,9).Value = 0  # I46
$ws.Cells.Item(46,10).Value = 1999  # J46
$ws.Cells.Item(46,11).Value = 0  # K46
$ws.Cells.Item(46,12).Value = 5997  # L46
$ws.Cells.Item(46,13).ClearContents()  # M46 removed
$ws.Cells.Item(46,14).Value = -6179  # N46
# row 92
$ws.Cells.Item(92,8).Value = 309.4  # H92
$ws.Cells.Item(92,9).Value = 448.5  # I92
$ws.Cells.Item(92,10).Value = 216.66667  # J92
$ws.Cells.Item(92,11).Value = 1345.5  # K92
$ws.Cells.Item(92,12).Value = 650.00001  # L92
$ws.Cells.Item(92,13).Value = -97.5  # M92
$ws.Cells.Item(92,14).Value = -3146.00001  # N92
# row 97
$ws.Cells.Item(97,8).Value = 405  # H97
$ws.Cells.Item(97,10).Value = 0  # J97
$ws.Cells.Item(97,12).Value = 0  # L97
$ws.Cells.Item(97,14).ClearContents()  # N97 removed
# row 109
$ws.Cells.Item(109,8).Value = 921.5  # H109
$ws.Cells.Item(109,9).Value = 921.5  # I109
$ws.Cells.Item(109,10).Value = 0  # J109
$ws.Cells.Item(109,11).Value = 2764.5  # K109
$ws.Cells.Item(109,12).Value = 0  # L109
$ws.Cells.Item(109,13).Value = -1724.5  # M109
$ws.Cells.Item(109,14).ClearContents()  # N109 removed
# row 117
$ws.Cells.Item(117,8).Value = 2163.7  # H117
$ws.Cells.Item(117,9).Value = 1862.6666  # I117
$ws.Cells.Item(117,10).Value = 2292.7144  # J117
$ws.Cells.Item(117,11).Value = 5587.9998  # K117
$ws.Cells.Item(117,12).Value = 6878.1432  # L117
$ws.Cells.Item(117,13).Value = -2145.9998  # M117
$ws.Cells.Item(117,14).Value = -13762.1432  # N117
# row 121
$ws.Cells.Item(121,8).Value = 2785.7144  # H121
# row 131
$ws.Cells.Item(131,8).Value = 1761.1852  # H131
$ws.Cells.Item(131,9).Value = 710.4167  # I131
$ws.Cells.Item(131,11).Value = 2131.2501  # K131
$ws.Cells.Item(131,13).Value = 2908.7499  # M131

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
# row 9
$ws.Cells.Item(9,8).Value = 334.5  # H9
$ws.Cells.Item(9,9).Value = 251.75  # I9
$ws.Cells.Item(9,11).Value = 251.75  # K9
$ws.Cells.Item(9,13).Value = -81.75  # M9
# row 11
$ws.Cells.Item(11,8).Value = 24574536  # H11
$ws.Cells.Item(11,9).Value = 33160104  # I11
$ws.Cells.Item(11,10).Value = 6449443  # J11
$ws.Cells.Item(11,11).Value = 33160104  # K11
$ws.Cells.Item(11,12).Value = 6449443  # L11
$ws.Cells.Item(11,13).Value = -33159965  # M11
$ws.Cells.Item(11,14).Value = -6449721  # N11
# row 12
$ws.Cells.Item(12,8).Value = 0  # H12
$ws.Cells.Item(12,9).Value = 0  # I12
$ws.Cells.Item(12,10).Value = 0  # J12
$ws.Cells.Item(12,11).Value = 0  # K12
$ws.Cells.Item(12,12).Value = 0  # L12
$ws.Cells.Item(12,13).ClearContents()  # M12 removed
$ws.Cells.Item(12,14).ClearContents()  # N12 removed
# row 31
$ws.Cells.Item(31,8).Value = 896.3333  # H31
$ws.Cells.Item(31,9).Value = 896.3333  # I31
$ws.Cells.Item(31,11).Value = 896.3333  # K31
$ws.Cells.Item(31,13).Value = -604.3333  # M31
# row 33
$ws.Cells.Item(33,8).Value = 3966.6667  # H33
$ws.Cells.Item(33,10).Value = 3966.6667  # J33
$ws.Cells.Item(33,12).Value = 3966.6667  # L33
$ws.Cells.Item(33,14).Value = -4470.6667  # N33
# row 36
$ws.Cells.Item(36,8).Value = 24067  # H36
$ws.Cells.Item(36,9).Value = 24067  # I36
$ws.Cells.Item(36,11).Value = 24067  # K36
$ws.Cells.Item(36,13).Value = -23582  # M36
# row 37
$ws.Cells.Item(37,8).Value = 896.3333  # H37
$ws.Cells.Item(37,9).Value = 896.3333  # I37
$ws.Cells.Item(37,11).Value = 896.3333  # K37
$ws.Cells.Item(37,13).Value = -619.3333  # M37
# row 40
$ws.Cells.Item(40,8).Value = 17466.666  # H40
$ws.Cells.Item(40,10).Value = 17466.666  # J40
$ws.Cells.Item(40,12).Value = 17466.666  # L40
$ws.Cells.Item(40,14).Value = -17768.666  # N40
# row 63
$ws.Cells.Item(63,8).Value = 41330.168  # H63
$ws.Cells.Item(63,9).Value = 24989  # I63
$ws.Cells.Item(63,10).Value = 44598.4  # J63
$ws.Cells.Item(63,11).Value = 24989  # K63
$ws.Cells.Item(63,12).Value = 44598.4  # L63
$ws.Cells.Item(63,13).Value = -24303  # M63
$ws.Cells.Item(63,14).Value = -45970.4  # N63
# row 66
$ws.Cells.Item(66,8).Value = 41330.168  # H66
$ws.Cells.Item(66,9).Value = 24989  # I66
$ws.Cells.Item(66,10).Value = 44598.4  # J66
$ws.Cells.Item(66,11).Value = 74967  # K66
$ws.Cells.Item(66,12).Value = 133795.2  # L66
$ws.Cells.Item(66,13).Value = -71535  # M66
$ws.Cells.Item(66,14).Value = -140659.2  # N66
# row 97
$ws.Cells.Item(97,8).Value = 112.55556  # H97
$ws.Cells.Item(97,9).Value = 157.6  # I97
$ws.Cells.Item(97,11).Value = 157.6  # K97
$ws.Cells.Item(97,13).Value = 338.4  # M97
# row 113
$ws.Cells.Item(113,8).Value = 4666.9  # H113
$ws.Cells.Item(113,9).Value = 2478.3333  # I113
$ws.Cells.Item(113,10).Value = 7949.75  # J113
$ws.Cells.Item(113,11).Value = 2478.3333  # K113
$ws.Cells.Item(113,12).Value = 7949.75  # L113
$ws.Cells.Item(113,13).Value = -308.3332999999998  # M113
$ws.Cells.Item(113,14).Value = -12289.75  # N113
# row 122
$ws.Cells.Item(122,8).Value = 456325.9  # H122
$ws.Cells.Item(122,9).Value = 626011.9  # I122
$ws.Cells.Item(122,10).Value = 3830  # J122
$ws.Cells.Item(122,11).Value = 1878035.7  # K122
$ws.Cells.Item(122,12).Value = 11490  # L122
$ws.Cells.Item(122,13).Value = -1875585.7  # M122
$ws.Cells.Item(122,14).Value = -16390  # N122
# row 136
$ws.Cells.Item(136,8).Value = 19999  # H136
$ws.Cells.Item(136,10).Value = 19999  # J136
$ws.Cells.Item(136,12).Value = 59997  # L136
$ws.Cells.Item(136,14).Value = -65097  # N136

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
# row 2
$ws.Cells.Item(2,8).Value = 58999.5  # H2
$ws.Cells.Item(2,9).Value = 9999  # I2
$ws.Cells.Item(2,11).Value = 9999  # K2
$ws.Cells.Item(2,13).Value = -9887  # M2
# row 7
$ws.Cells.Item(7,8).Value = 7396.75  # H7
$ws.Cells.Item(7,9).Value = 7183.4287  # I7
$ws.Cells.Item(7,10).Value = 8890  # J7
$ws.Cells.Item(7,11).Value = 7183.4287  # K7
$ws.Cells.Item(7,12).Value = 8890  # L7
$ws.Cells.Item(7,13).Value = -7071.4287  # M7
$ws.Cells.Item(7,14).Value = -9114  # N7
# row 14
$ws.Cells.Item(14,8).Value = 14856.286  # H14
$ws.Cells.Item(14,9).Value = 14832.366  # I14
$ws.Cells.Item(14,10).Value = 14999.8  # J14
$ws.Cells.Item(14,11).Value = 14832.366  # K14
$ws.Cells.Item(14,12).Value = 14999.8  # L14
$ws.Cells.Item(14,13).Value = -14660.366  # M14
$ws.Cells.Item(14,14).Value = -15343.8  # N14
# row 16
$ws.Cells.Item(16,8).Value = 841  # H16
$ws.Cells.Item(16,9).Value = 793.3333  # I16
$ws.Cells.Item(16,10).Value = 912.5  # J16
$ws.Cells.Item(16,11).Value = 793.3333  # K16
$ws.Cells.Item(16,12).Value = 912.5  # L16
$ws.Cells.Item(16,13).Value = -623.3333  # M16
$ws.Cells.Item(16,14).Value = -1252.5  # N16
# row 20
$ws.Cells.Item(20,8).Value = 20000  # H20
$ws.Cells.Item(20,9).Value = 0  # I20
$ws.Cells.Item(20,10).Value = 20000  # J20
$ws.Cells.Item(20,11).Value = 0  # K20
$ws.Cells.Item(20,12).Value = 20000  # L20
$ws.Cells.Item(20,13).ClearContents()  # M20 removed
$ws.Cells.Item(20,14).Value = -20452  # N20
# row 22
$ws.Cells.Item(22,8).Value = 1862.5  # H22
$ws.Cells.Item(22,9).Value = 775  # I22
$ws.Cells.Item(22,10).Value = 2950  # J22
$ws.Cells.Item(22,11).Value = 775  # K22
$ws.Cells.Item(22,12).Value = 2950  # L22
$ws.Cells.Item(22,13).Value = -480  # M22
$ws.Cells.Item(22,14).Value = -3540  # N22
# row 27
$ws.Cells.Item(27,8).Value = 1862.5  # H27
$ws.Cells.Item(27,9).Value = 775  # I27
$ws.Cells.Item(27,10).Value = 2950  # J27
$ws.Cells.Item(27,11).Value = 775  # K27
$ws.Cells.Item(27,12).Value = 2950  # L27
$ws.Cells.Item(27,13).Value = -668  # M27
$ws.Cells.Item(27,14).Value = -3164  # N27
# row 46
$ws.Cells.Item(46,8).Value = 4304.6  # H46
$ws.Cells.Item(46,9).Value = 2174.3333  # I46
$ws.Cells.Item(46,10).Value = 7500  # J46
$ws.Cells.Item(46,11).Value = 2174.3333  # K46
$ws.Cells.Item(46,12).Value = 7500  # L46
$ws.Cells.Item(46,13).Value = -1986.3333  # M46
$ws.Cells.Item(46,14).Value = -7876  # N46
# row 55
$ws.Cells.Item(55,8).Value = 998.17645  # H55
$ws.Cells.Item(55,9).Value = 917.2  # I55
$ws.Cells.Item(55,10).Value = 1113.8572  # J55
$ws.Cells.Item(55,11).Value = 917.2  # K55
$ws.Cells.Item(55,12).Value = 1113.8572  # L55
$ws.Cells.Item(55,13).Value = -744.2  # M55
$ws.Cells.Item(55,14).Value = -1459.8572  # N55
# row 61
$ws.Cells.Item(61,8).Value = 5220.231  # H61
$ws.Cells.Item(61,9).Value = 4152.3335  # I61
$ws.Cells.Item(61,11).Value = 4152.3335  # K61
$ws.Cells.Item(61,13).Value = -3950.3335  # M61
# row 93
$ws.Cells.Item(93,8).Value = 9642.857  # H93
$ws.Cells.Item(93,9).Value = 16250  # I93
$ws.Cells.Item(93,10).Value = 7000  # J93
$ws.Cells.Item(93,11).Value = 16250  # K93
$ws.Cells.Item(93,12).Value = 7000  # L93
$ws.Cells.Item(93,13).Value = -15002  # M93
$ws.Cells.Item(93,14).Value = -9496  # N93
# row 113
$ws.Cells.Item(113,8).Value = 5220.231  # H113
$ws.Cells.Item(113,9).Value = 4152.3335  # I113
$ws.Cells.Item(113,11).Value = 4152.3335  # K113
$ws.Cells.Item(113,13).Value = -1982.3335  # M113
# row 126
$ws.Cells.Item(126,8).Value = 7396.75  # H126
$ws.Cells.Item(126,9).Value = 7183.4287  # I126
$ws.Cells.Item(126,10).Value = 8890  # J126
$ws.Cells.Item(126,11).Value = 21550.2861  # K126
$ws.Cells.Item(126,12).Value = 26670  # L126
$ws.Cells.Item(126,13).Value = -19080.2861  # M126
$ws.Cells.Item(126,14).Value = -31610  # N126

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
# row 4
$ws.Cells.Item(4,8).Value = 31557.154  # H4
$ws.Cells.Item(4,10).Value = 0  # J4
$ws.Cells.Item(4,12).Value = 0  # L4
$ws.Cells.Item(4,14).ClearContents()  # N4 removed
# row 5
$ws.Cells.Item(5,8).Value = 20000000  # H5
$ws.Cells.Item(5,9).Value = 20000000  # I5
$ws.Cells.Item(5,11).Value = 20000000  # K5
$ws.Cells.Item(5,13).Value = -19999888  # M5
# row 23
$ws.Cells.Item(23,8).Value = 309.2  # H23
$ws.Cells.Item(23,9).Value = 123  # I23
$ws.Cells.Item(23,10).Value = 355.75  # J23
$ws.Cells.Item(23,11).Value = 123  # K23
$ws.Cells.Item(23,12).Value = 355.75  # L23
$ws.Cells.Item(23,13).Value = 106  # M23
$ws.Cells.Item(23,14).Value = -813.75  # N23
# row 32
$ws.Cells.Item(32,8).Value = 4499  # H32
$ws.Cells.Item(32,9).Value = 4499  # I32
$ws.Cells.Item(32,11).Value = 4499  # K32
$ws.Cells.Item(32,13).Value = -4182  # M32
# row 41
$ws.Cells.Item(41,8).Value = 36973.4  # H41
$ws.Cells.Item(41,9).Value = 40000  # I41
$ws.Cells.Item(41,11).Value = 40000  # K41
$ws.Cells.Item(41,13).Value = -39610  # M41
# row 58
$ws.Cells.Item(58,8).Value = 0  # H58
$ws.Cells.Item(58,9).Value = 0  # I58
$ws.Cells.Item(58,11).Value = 0  # K58
$ws.Cells.Item(58,13).ClearContents()  # M58 removed
# row 64
$ws.Cells.Item(64,8).Value = 90000  # H64
$ws.Cells.Item(64,9).Value = 90000  # I64
$ws.Cells.Item(64,11).Value = 90000  # K64
$ws.Cells.Item(64,13).Value = -89752  # M64
# row 67
$ws.Cells.Item(67,8).Value = 90000  # H67
$ws.Cells.Item(67,9).Value = 90000  # I67
$ws.Cells.Item(67,11).Value = 90000  # K67
$ws.Cells.Item(67,13).Value = -89142  # M67
# row 70
$ws.Cells.Item(70,8).Value = 90000  # H70
$ws.Cells.Item(70,9).Value = 90000  # I70
$ws.Cells.Item(70,11).Value = 90000  # K70
$ws.Cells.Item(70,13).Value = -89685  # M70
# row 73
$ws.Cells.Item(73,8).Value = 90000  # H73
$ws.Cells.Item(73,9).Value = 90000  # I73
$ws.Cells.Item(73,11).Value = 90000  # K73
$ws.Cells.Item(73,13).Value = -88908  # M73
# row 76
$ws.Cells.Item(76,8).Value = 30000  # H76
$ws.Cells.Item(76,10).Value = 30000  # J76
$ws.Cells.Item(76,12).Value = 30000  # L76
$ws.Cells.Item(76,14).Value = -30630  # N76
# row 79
$ws.Cells.Item(79,8).Value = 30000  # H79
$ws.Cells.Item(79,10).Value = 30000  # J79
$ws.Cells.Item(79,12).Value = 30000  # L79
$ws.Cells.Item(79,14).Value = -32184  # N79
# row 95
$ws.Cells.Item(95,8).Value = 10000  # H95
$ws.Cells.Item(95,10).Value = 10000  # J95
$ws.Cells.Item(95,12).Value = 10000  # L95
$ws.Cells.Item(95,14).Value = -15492  # N95
# row 113
$ws.Cells.Item(113,8).Value = 573.3333  # H113
$ws.Cells.Item(113,9).Value = 428.66666  # I113
$ws.Cells.Item(113,10).Value = 718  # J113
$ws.Cells.Item(113,11).Value = 1285.99998  # K113
$ws.Cells.Item(113,12).Value = 2154  # L113
$ws.Cells.Item(113,13).Value = 884.0000199999999  # M113
$ws.Cells.Item(113,14).Value = -6494  # N113
# row 122
$ws.Cells.Item(122,8).Value = 3051.6667  # H122
$ws.Cells.Item(122,9).Value = 3101.8635  # I122
$ws.Cells.Item(122,10).Value = 2499.5  # J122
$ws.Cells.Item(122,11).Value = 9305.5905  # K122
$ws.Cells.Item(122,12).Value = 7498.5  # L122
$ws.Cells.Item(122,13).Value = -6855.5905  # M122
$ws.Cells.Item(122,14).Value = -12398.5  # N122
# row 126
$ws.Cells.Item(126,8).Value = 5312.55  # H126
$ws.Cells.Item(126,9).Value = 3371.2222  # I126
$ws.Cells.Item(126,10).Value = 6900.909  # J126
$ws.Cells.Item(126,11).Value = 10113.6666  # K126
$ws.Cells.Item(126,12).Value = 20702.727  # L126
$ws.Cells.Item(126,13).Value = -7643.6666  # M126
$ws.Cells.Item(126,14).Value = -25642.727  # N126
# row 131
$ws.Cells.Item(131,8).Value = 73571.664  # H131
$ws.Cells.Item(131,10).Value = 73571.664  # J131
$ws.Cells.Item(131,12).Value = 73571.664  # L131
$ws.Cells.Item(131,14).Value = -83651.664  # N131
# row 132
$ws.Cells.Item(132,8).Value = 1932.9286  # H132
$ws.Cells.Item(132,9).Value = 1773.9231  # I132
$ws.Cells.Item(132,11).Value = 5321.7693  # K132
$ws.Cells.Item(132,13).Value = -2791.7693  # M132
# row 136
$ws.Cells.Item(136,8).Value = 3723.56  # H136
$ws.Cells.Item(136,9).Value = 2907.3635  # I136
$ws.Cells.Item(136,11).Value = 8722.0905  # K136
$ws.Cells.Item(136,13).Value = -6172.0905  # M136
